# blade of arcana years
#
# The "year" column (A) for the Super Scenario Supplement block (and a
# couple of neighboring rows) was never filled in. This fills in the
# correct year for every row in that block, re-sorts the block by year
# (ascending, stable) so the rows land in chronological order like the
# rest of the sheet, and then - matching the sheet's existing
# convention of leaving the year cell blank when it repeats the value
# directly above it - clears the year back out for the rows where it
# would just repeat the previous row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the real (previously missing) publication year for rows 6-30.
$ws.Range("A6").Value  = 2003   # Chains of Fortune
$ws.Range("A7").Value  = 2004   # Melody of Minstrel
$ws.Range("A8").Value  = 2001   # Road of Glory: The 2nd Edition
$ws.Range("A9").Value  = 2001   # Super Scenario Supplement Vol. 1
$ws.Range("A10").Value = 2002   # Super Scenario Supplement Vol. 2
$ws.Range("A11").Value = 2002   # Super Scenario Supplement Vol. 3
$ws.Range("A12").Value = 2002   # Super Scenario Supplement Vol. 4
$ws.Range("A13").Value = 2002   # Super Scenario Supplement Vol. 5
$ws.Range("A14").Value = 2002   # Super Scenario Supplement Vol. 6
$ws.Range("A15").Value = 2002   # Super Scenario Supplement Vol. 7
$ws.Range("A16").Value = 2003   # Super Scenario Supplement Vol. 8
$ws.Range("A17").Value = 2003   # Super Scenario Supplement Vol. 9
$ws.Range("A18").Value = 2003   # Super Scenario Supplement Vol. 10
$ws.Range("A19").Value = 2003   # Super Scenario Supplement Vol. 11
$ws.Range("A20").Value = 2003   # Super Scenario Supplement Vol. 12
$ws.Range("A21").Value = 2003   # Super Scenario Supplement Vol. 13
$ws.Range("A22").Value = 2003   # Super Scenario Supplement Vol. 14
$ws.Range("A23").Value = 2004   # Super Scenario Supplement Vol. 15
$ws.Range("A24").Value = 2004   # Super Scenario Supplement Vol. 16
$ws.Range("A25").Value = 2004   # Super Scenario Supplement Vol. 17
$ws.Range("A26").Value = 2004   # Super Scenario Supplement Vol. 18
$ws.Range("A27").Value = 2004   # Super Scenario Supplement Vol. 19
$ws.Range("A28").Value = 2005   # Super Scenario Supplement Vol. 20
$ws.Range("A29").Value = 2005   # Super Scenario Supplement Vol. 21
$ws.Range("A30").Value = 2005   # Super Scenario Supplement Vol. 22

# Re-sort the whole table chronologically by year (rows 2-44, matching
# the range the author's last Data > Sort used).
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A44"))
$sortObj.SetRange($ws.Range("A2:F44"))
$sortObj.Header = -4163
$sortObj.Apply()

# Match the sheet's convention of blanking a year cell when it is the
# same as the cell directly above it.
$ws.Range("A6").ClearContents()
$ws.Range("A13").ClearContents()
$ws.Range("A21").ClearContents()

# Update the selection left behind from editing to reflect where the
# user ended up after the sort (row 21).
$ws.Rows(21).Select()
